$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 341, which shifts existing rows 341:445 down to 342:446
$ws.Rows.Item(341).Insert()

# Populate the newly inserted row 341 with the new data
$ws.Cells.Item(341, 1).Value = 9
$ws.Cells.Item(341, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(341, 3).Value = "Metropolitana"
$ws.Cells.Item(341, 4).Value = 45093
$ws.Cells.Item(341, 5).Value = 13
$ws.Cells.Item(341, 6).Value = 300000001
$ws.Cells.Item(341, 7).Value = "Rabanito"
$ws.Cells.Item(341, 8).Value = "Sin especificar"
$ws.Cells.Item(341, 9).Value = "Primera"
$ws.Cells.Item(341, 10).Value = 7000
$ws.Cells.Item(341, 11).Value = 3000
$ws.Cells.Item(341, 12).Value = 3000
$ws.Cells.Item(341, 13).Value = 3000
$ws.Cells.Item(341, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(341, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(341, 16).Value = 30
$ws.Cells.Item(341, 17).Value = 100
$ws.Cells.Item(341, 18).Value = "Hortaliza"
